# Insert a new data row at row 329 (pushing the existing rows 329-349 down
# to 330-350) and populate it with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(329).Insert()

$ws.Range("A329").Value2 = 3
$ws.Range("B329").Value2 = "Femacal de La Calera"
$ws.Range("C329").Value2 = "Coquimbo"
$ws.Range("D329").Value2 = 44610
$ws.Range("E329").Value2 = 5
$ws.Range("F329").Value2 = "Fruta"
$ws.Range("G329").Value2 = 100108
$ws.Range("H329").Value2 = "Tropicales y subtropicales"
$ws.Range("I329").Value2 = 100108002
$ws.Range("J329").Value2 = "Mango"
$ws.Range("K329").Value2 = "Sin especificar"
$ws.Range("L329").Value2 = "Primera"
$ws.Range("M329").Value2 = 456
$ws.Range("N329").Value2 = 7000
$ws.Range("O329").Value2 = 7000
$ws.Range("P329").Value2 = 7000
$ws.Range("Q329").Value2 = "`$/bandeja 4 kilos"
$ws.Range("R329").Value2 = "Perú"
$ws.Range("S329").Value2 = 1750
$ws.Range("T329").Value2 = 4
